# PT-11, PT-15, PT-16, PT-17, PT-18, PT-19
# Restructure the "Kapazität" capacity sheet: split the former single
# "Summe intern" total row into separate internal/external name lists with
# their own sum rows, and refresh the rolling month header.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item('Kapazität')

# --- Make room for two new "external" resource rows at 7 and 8 ---
# (this pushes the old row 7 "Summe intern" row down to row 9, and all
# of its SUM(B2:B6)-style formulas shift automatically with it)
$ws.Range("A7:A8").EntireRow.Insert()

# --- New "Summe extern" row at 11 (row 10 stays blank as a spacer) ---
$ws.Range("A11").Value = "Summe extern"
$ws.Range("B11:M11").Formula = "=SUM(B7:B8)"

# The old "Summe intern" row covered 21 months (B:V); the rebuilt sheet only
# tracks 12 (B:M), so drop the trailing SUM(..)-formulas that rode along
# when the row shifted from 7 to 9.
$ws.Range("N9:V9").Clear()

# --- Rename the five internal resource placeholders (A2:A6) and fill
#     in the two new external ones (A7:A8) ---
$ws.Range("A2").Value = "<Name 1>"
$ws.Range("A3").Value = "<Name 2>"
$ws.Range("A4").Value = "<Name 3>"
$ws.Range("A5").Value = "<Name 4>"
$ws.Range("A6").Value = "<Name 5>"
$ws.Range("A7").Value = "<Name Extern1>"
$ws.Range("A8").Value = "<Name Extern2>"

# Copy the "Summe intern" row's formatting (fill/alignment) across A11:JU11
# so the new total row matches the look of the existing one.
$ws.Range("A9:V9").Copy()
$ws.Range("A11:JU11").PasteSpecial(-4122)
$ws.Rows.Item(11).RowHeight = 30.75

# --- Refresh the rolling 12-month header row (shift one year forward) ---
$ws.Range("B1").Value = 41821
$ws.Range("C1").Value = 41852
$ws.Range("D1").Value = 41883
$ws.Range("E1").Value = 41913
$ws.Range("F1").Value = 41944
$ws.Range("G1").Value = 41974
$ws.Range("H1").Value = 42005
$ws.Range("I1").Value = 42036
$ws.Range("J1").Value = 42064
$ws.Range("K1").Value = 42095
$ws.Range("L1").Value = 42125
$ws.Range("M1").Value = 42156
$ws.Range("N1:V1").ClearContents()
$ws.Range("BD1:BL1").Clear()

$ws.Range("M11").Select()

# --- Defined names: repoint intern_sum, add extern_sum ---
$n = $wb.Names.Item('intern_sum')
$n.RefersTo = '=Kapazität!$A$9'
$wb.Names.Add('extern_sum', '=Kapazität!$A$11')

# --- Summary sheet: "Produkt 4" rows used to reuse the old shared string
#     slot that "Tom " occupied; after the rename that slot holds "Produkt 4"
#     itself, so repoint B8/B18/B28/B37/B47 there. ---
$wsSummary = $wb.Worksheets.Item('Summary')
$wsSummary.Range("B8").Value = "Produkt 4"
$wsSummary.Range("B18").Value = "Produkt 4"
$wsSummary.Range("B28").Value = "Produkt 4"
$wsSummary.Range("B37").Value = "Produkt 4"
$wsSummary.Range("B47").Value = "Produkt 4"
